$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): re-label + move the "centered" styling ----
# Clear all current formatting off the header cells first so no stray
# style survives the column reshuffle, then re-apply it where it belongs.
$ws.Range("A1:J1").ClearFormats()

$ws.Range("A1").Value = "Subject"
$ws.Range("B1").Value = "way_of_speech"
$ws.Range("C1").Value = "socio_economic"
$ws.Range("D1").Value = "ethnicity_skin_color"
$ws.Range("E1").Value = "personality"
$ws.Range("F1").Value = "dress_propeties"
$ws.Range("G1").Value = "political_affiliation"
$ws.Range("H1").Value = "hobbies"
$ws.Range("I1").Value = "body_size"
$ws.Range("J1").Value = "intelligence"

# Re-apply the bold font on A1 ("Subject")
$ws.Range("A1").Font.Bold = $true

# Re-apply the vertical-center alignment on D1 ("ethnicity_skin_color")
# and H1 ("hobbies")
$ws.Range("D1").VerticalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4108

# ---- Data rows ----
# Wipe the old 2-9 block, then write the fresh 4-row dataset (rows 2-5)
$ws.Range("A2:J9").Clear()

$data = @(
    @(99881, 62, 71, 61, 44, 50, 88, 68, 14, 81),
    @(99882, 66, 75, 20, 38, 19, 29, 18, 14, 67),
    @(99883, 49, 85, 11, 66, 94, 50, 55, 64, 59),
    @(99884, 59, 65, 30, 12, 41, 43, 15, 59, 69)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# ---- Selection ----
$ws.Range("A6:XFD8").Select() | Out-Null
